$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8095196980875699
$ws.Range("C2").Value = 0.2248015017041496
$ws.Range("D2").Value = 0.01752849954091573
$ws.Range("E2").Value = 0.1154778989324541
$ws.Range("F2").Value = 0.4666915567120284
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.3178451623716221
$ws.Range("M2").Value = 0.3456511347932647
$ws.Range("O2").Value = 1.49006755305814
$ws.Range("B3").Value = 0.7078771821256851
$ws.Range("C3").Value = 0.1978619620783775
$ws.Range("D3").Value = 0.0156586582833711
$ws.Range("E3").Value = 0.1108315602518104
$ws.Range("F3").Value = 0.463561768233653
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.3245545926818263
$ws.Range("M3").Value = 0.3069021009581689
$ws.Range("O3").Value = 1.493432419653843
$ws.Range("B4").Value = 0.6452677079924456
$ws.Range("C4").Value = 0.1812406018929664
$ws.Range("D4").Value = 0.01450505921919643
$ws.Range("E4").Value = 0.1081095712119478
$ws.Range("F4").Value = 0.4621030246409745
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.3290083440213127
$ws.Range("M4").Value = 0.2831425954706575
$ws.Range("O4").Value = 1.497032555072366
$ws.Range("B5").Value = 0.6197049179057217
$ws.Range("C5").Value = 0.1744474129770595
$ws.Range("D5").Value = 0.01403360736040327
$ws.Range("E5").Value = 0.1070329263859655
$ws.Range("F5").Value = 0.4616246226559753
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.330907036048556
$ws.Range("M5").Value = 0.2734686878302597
$ws.Range("O5").Value = 1.4988840719856
$ws.Range("B6").Value = 0.6154573261830762
$ws.Range("C6").Value = 0.1733182209742665
$ws.Range("D6").Value = 0.01395524256825809
$ws.Range("E6").Value = 0.106856109101674
$ws.Range("F6").Value = 0.4615521826102693
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.3312273629607105
$ws.Range("M6").Value = 0.2718628461408557
$ws.Range("O6").Value = 1.49921469478403
$ws.Range("B7").Value = 0.6449231553299626
$ws.Range("C7").Value = 0.1811490664583459
$ws.Range("D7").Value = 0.01449870647803664
$ws.Range("E7").Value = 0.108094919683964
$ws.Range("F7").Value = 0.4620961033594568
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.3290336116940811
$ws.Range("M7").Value = 0.2830120960950069
$ws.Range("O7").Value = 1.497055970491829
$ws.Range("B8").Value = 0.7745158556203933
$ws.Range("C8").Value = 0.2155296374935745
$ws.Range("D8").Value = 0.01688494164645249
$ws.Range("E8").Value = 0.1138484804919742
$ws.Range("F8").Value = 0.4655160954658299
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.320089063765387
$ws.Range("M8").Value = 0.3322837554055909
$ws.Range("O8").Value = 1.490908555598367
$ws.Range("B9").Value = 1.027001681523927
$ws.Range("C9").Value = 0.2822995273424738
$ws.Range("D9").Value = 0.02151934725362281
$ws.Range("E9").Value = 0.1261849593141164
$ws.Range("F9").Value = 0.475913767220014
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.3052118937926558
$ws.Range("M9").Value = 0.4291645336342214
$ws.Range("O9").Value = 1.491087815790308
$ws.Range("B10").Value = 1.211446582613291
$ws.Range("C10").Value = 0.3309468636891779
$ws.Range("D10").Value = 0.02489542674197054
$ws.Range("E10").Value = 0.1359131060903138
$ws.Range("F10").Value = 0.4858292803131619
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.2959211570312767
$ws.Range("M10").Value = 0.5005096146014836
$ws.Range("O10").Value = 1.498766673317988
$ws.Range("B11").Value = 1.295115976528677
$ws.Range("C11").Value = 0.3529868942402743
$ws.Range("D11").Value = 0.02642476604707866
$ws.Range("E11").Value = 0.140487820143818
$ws.Range("F11").Value = 0.4908400923093552
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.2920539943756069
$ws.Range("M11").Value = 0.533004950474151
$ws.Range("O11").Value = 1.503917999767623
$ws.Range("B12").Value = 1.326764311612806
$ws.Range("C12").Value = 0.3613196692631391
$ws.Range("D12").Value = 0.02700292933707971
$ws.Range("E12").Value = 0.1422419899097065
$ws.Range("F12").Value = 0.4928099070009324
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.2906415471166284
$ws.Range("M12").Value = 0.5453158901342334
$ws.Range("O12").Value = 1.506108606981769
$ws.Range("B13").Value = 1.31994987797151
$ws.Range("C13").Value = 0.3595256548818497
$ws.Range("D13").Value = 0.02687845498419961
$ws.Range("E13").Value = 0.1418632217955249
$ws.Range("F13").Value = 0.4923824489407167
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.2909434269713778
$ws.Range("M13").Value = 0.5426642588557087
$ws.Range("O13").Value = 1.505626127055791
$ws.Range("B14").Value = 1.297720423541534
$ws.Range("C14").Value = 0.3536727054982691
$ws.Range("D14").Value = 0.02647235141961346
$ws.Range("E14").Value = 0.1406316973284447
$ws.Range("F14").Value = 0.4910006982770909
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.2919367483219446
$ws.Range("M14").Value = 0.5340176660181726
$ws.Range("O14").Value = 1.504093405760898
$ws.Range("B15").Value = 1.284099577754091
$ws.Range("C15").Value = 0.3500858609158968
$ws.Range("D15").Value = 0.02622347461790753
$ws.Range("E15").Value = 0.1398802057122737
$ws.Range("F15").Value = 0.4901637673437875
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.2925519625321265
$ws.Range("M15").Value = 0.5287221124437309
$ws.Range("O15").Value = 1.503185856789116
$ws.Range("B16").Value = 1.205973716162532
$ws.Range("C16").Value = 0.329504654474249
$ws.Range("D16").Value = 0.0247953478279328
$ws.Range("E16").Value = 0.1356171726912194
$ws.Range("F16").Value = 0.485511911534509
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.2961811401197103
$ws.Range("M16").Value = 0.498386765820257
$ws.Range("O16").Value = 1.498463514346895
$ws.Range("B17").Value = 1.157984578316018
$ws.Range("C17").Value = 0.3168554566060777
$ws.Range("D17").Value = 0.0239175590395746
$ws.Range("E17").Value = 0.1330404368752909
$ws.Range("F17").Value = 0.4827865469423429
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.2984997502116009
$ws.Range("M17").Value = 0.4797872078058987
$ws.Range("O17").Value = 1.495992254748501
$ws.Range("B18").Value = 1.13036042098031
$ws.Range("C18").Value = 0.309571534463231
$ws.Range("D18").Value = 0.02341207258198352
$ws.Range("E18").Value = 0.1315724012997066
$ws.Range("F18").Value = 0.4812660461846079
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.2998671485818178
$ws.Range("M18").Value = 0.4690930144836045
$ws.Range("O18").Value = 1.494726839439124
$ws.Range("B19").Value = 1.121003613107973
$ws.Range("C19").Value = 0.3071038859982309
$ws.Range("D19").Value = 0.02324082064728827
$ws.Range("E19").Value = 0.131077749469803
$ws.Range("F19").Value = 0.4807593005557962
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.3003359212279939
$ws.Range("M19").Value = 0.4654727969526391
$ws.Range("O19").Value = 1.494325137825598
$ws.Range("B20").Value = 1.163095400105192
$ws.Range("C20").Value = 0.3182028613953491
$ws.Range("D20").Value = 0.02401106408915865
$ws.Range("E20").Value = 0.1333132797221666
$ws.Range("F20").Value = 0.483071793531046
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.2982494303234553
$ws.Range("M20").Value = 0.4817667720185739
$ws.Range("O20").Value = 1.496239169741216
$ws.Range("B21").Value = 1.304250727739316
$ws.Range("C21").Value = 0.3553922227492308
$ws.Range("D21").Value = 0.02659166033619442
$ws.Range("E21").Value = 0.1409928307931452
$ws.Range("F21").Value = 0.4914045857643288
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.2916435730649383
$ws.Range("M21").Value = 0.5365572271058738
$ws.Range("O21").Value = 1.504537080346552
$ws.Range("B22").Value = 1.396296551371961
$ws.Range("C22").Value = 0.3796198878458767
$ws.Range("D22").Value = 0.02827259199222709
$ws.Range("E22").Value = 0.1461392369728998
$ws.Range("F22").Value = 0.497272320586049
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.2876292866373369
$ws.Range("M22").Value = 0.5723989087631622
$ws.Range("O22").Value = 1.511359311726665
$ws.Range("B23").Value = 1.347189452966518
$ws.Range("C23").Value = 0.3666963641765619
$ws.Range("D23").Value = 0.02737597470795095
$ws.Range("E23").Value = 0.1433807319422655
$ws.Range("F23").Value = 0.4941018756241533
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.2897439595537605
$ws.Range("M23").Value = 0.5532665473143226
$ws.Range("O23").Value = 1.507589659789801
$ws.Range("B24").Value = 1.160784905117168
$ws.Range("C24").Value = 0.3175937361992567
$ws.Range("D24").Value = 0.02396879305113231
$ws.Range("E24").Value = 0.1331898858714737
$ws.Range("F24").Value = 0.4829426892079951
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.2983624928167501
$ws.Range("M24").Value = 0.4808718143371351
$ws.Range("O24").Value = 1.496127055689129
$ws.Range("B25").Value = 0.9588792958925865
$ws.Range("C25").Value = 0.2643073860541847
$ws.Range("D25").Value = 0.02027057965481305
$ws.Range("E25").Value = 0.1227325180173082
$ws.Range("F25").Value = 0.4727030551500206
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.3089498032192175
$ws.Range("M25").Value = 0.4029270774830707
$ws.Range("O25").Value = 1.489720359830613
